$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04) for every data row
# (rows 2 through 72).
$ws.Range("C2:C72").Value = 45203
